# "fix table grade filling"
#
# The Контесты sheet tracks, per student (rows 4-6), whether each of the
# contest problems A..O* (columns B..P) was solved (1) or not (0), with
# column Q holding the total Score for that student. The grade-filling
# logic previously left every problem/Score cell at 0 for all students;
# this fixes it to contain the real per-problem results and the resulting
# total score (sum of the problem columns) for each student.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row => values for columns B..Q (15 problem columns + Score)
$grades = @{
    4 = @(1, 1, 1, 1, 1, 0, 0, 0, 1, 1, 0, 0, 0, 0, 0, 7)   # Абабков Даниил
    5 = @(1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 0, 0, 1, 0, 0, 11)  # Гусев Александр Дмитриевич
    6 = @(1, 1, 1, 1, 1, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 5)   # hse-compds-2022-46
}

foreach ($row in $grades.Keys) {
    $values = $grades[$row]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = 2 + $i   # column B is index 2
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}
